$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fix the cached "datetimeFigureOut" date placeholder text on every slide
#    layout and on the slide master: "27.6.2022 г." -> "27.6.2022".
# ---------------------------------------------------------------------------
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = -1
        try {
            $phType = $shp.PlaceholderFormat.Type
        } catch {
            $phType = -1
        }
        if ($phType -eq $ppPlaceholderDate -and $shp.HasTextFrame) {
            $curText = $shp.TextFrame.TextRange.Text
            if ($curText -eq "27.6.2022 г." -or $curText -like "27.6.2022*") {
                if ($curText -ne "27.6.2022") {
                    $shp.TextFrame.TextRange.Text = "27.6.2022"
                }
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Remove the empty subtitle placeholder ("Подзаглавие 2") from slide 1.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "Подзаглавие 2") {
        $shp.Delete()
    }
}
